$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell A4 value from 5 to 7
$ws.Range("A4").Value = 7

# Update selection to C7
$ws.Range("C7").Select()

# Add data validations
$xlValidateCustom = 7
$xlValidAlertStop = 1

$range1 = $ws.Range("A2:A564")
$range1.Validation.Add($xlValidateCustom, $xlValidAlertStop, 1, 'COUNTIF($A:$A,A2)=1')
$range1.Validation.ErrorTitle = "Duplicate Scene ID"
$range1.Validation.ErrorMessage = "This scene ID already exists. You are either using the wrong ID, or you need to update an existing entry."
$range1.Validation.ShowInput = $true
$range1.Validation.ShowError = $true

$range2 = $ws.Range("B1:B1048576")
$range2.Validation.Add($xlValidateCustom, $xlValidAlertStop, 1, 'COUNTIF($B:$B,B1)=1')
$range2.Validation.ErrorTitle = "Duplicate Title"
$range2.Validation.ErrorMessage = "You have entered a duplicate Title. This scene either already exists, and you need to find it, or you need to create a more meaningful (and unique) title."
$range2.Validation.ShowInput = $true
$range2.Validation.ShowError = $true

# Set page orientation to portrait
$ws.PageSetup.Orientation = 1


